$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for rows 2-210 to new date serial 45184
$ws.Range("C2:C210").Value = 45184

# Append new row 211 with data for A 43077-2023
$ws.Range("A211").Value = "A 43077-2023"
$ws.Range("B211").Value = 45182
$ws.Range("C211").Value = 45184
$ws.Range("B211:C211").NumberFormat = "YYYY-MM-DD"
$ws.Range("D211").Value = "ÖREBRO LÄN"
$ws.Range("E211").Value = "LEKEBERG"
$ws.Range("G211").Value = 1
$ws.Range("H211").Value = 0
$ws.Range("I211").Value = 0
$ws.Range("J211").Value = 0
$ws.Range("K211").Value = 0
$ws.Range("L211").Value = 0
$ws.Range("M211").Value = 0
$ws.Range("N211").Value = 0
$ws.Range("O211").Value = 0
$ws.Range("P211").Value = 0
$ws.Range("Q211").Value = 0

# Copy the R210 cell's style (wrap text, empty) down to R211
$ws.Range("R210").Copy()
$ws.Range("R211").PasteSpecial(-4122)

# Row 210 picks up an explicit custom height as a side effect of the edit
$ws.Rows.Item(210).RowHeight = 15
